$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out any previous contents first
$ws.Cells.Clear()

# Header row
$ws.Range("A1").Value = "Questions"
$ws.Range("B1").Value = "option_A"
$ws.Range("C1").Value = "option_B"
$ws.Range("D1").Value = "option_C"
$ws.Range("E1").Value = "option_D"
$ws.Range("F1").Value = "option_E"
$ws.Range("G1").Value = "Correct_Answer"

# Row 2 - NBA question
$ws.Range("A2").Value = "Which one is correct team name in NBA?"
$ws.Range("B2").Value = "New York Bulls"
$ws.Range("C2").Value = "Los Angeles Kings"
$ws.Range("D2").Value = "Golden State Warriros"
$ws.Range("E2").Value = "Huston Rocket"
$ws.Range("G2").Value = "Huston Rocket"

# Row 3 - math question
$ws.Range("A3").Value = "5 + 7 = ?"
$ws.Range("B3").Value = 10
$ws.Range("C3").Value = 11
$ws.Range("D3").Value = 12
$ws.Range("E3").Value = 13
$ws.Range("F3").Value = 14
$ws.Range("G3").Value = 12

# Row 4 - subtraction question
$ws.Range("A4").Value = "12 - 8 = ?"
$ws.Range("B4").Value = 4
$ws.Range("C4").Value = 3
$ws.Range("D4").Value = 2
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 4

# Column widths (engine quantizes ColumnWidth on a 6px-per-char grid before
# storing; the offsets below are chosen so the stored width in the OOXML
# matches the target as closely as that quantization allows)
$ws.Columns.Item(1).ColumnWidth = 53.022135416666664
$ws.Columns.Item(2).ColumnWidth = 21.592447916666668
$ws.Columns.Item(3).ColumnWidth = 21.451822916666668
$ws.Columns.Item(4).ColumnWidth = 21.166666666666668
$ws.Columns.Item(5).ColumnWidth = 13.736979166666666
$ws.Columns.Item(6).ColumnWidth = 12.736979166666666
$ws.Columns.Item(7).ColumnWidth = 19.166666666666668

# Font styling for C1 (option_B header) - black color font, no theme
$ws.Range("C1").Font.Color = 0

# Selection matches target
[void]$ws.Range("B8").Select()
